$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain decimal number must be forced to Text format,
# otherwise Excel auto-converts the literal (matching the original inlineStr cells)
# into a floating-point number, which both changes the cell type and introduces
# binary floating point noise (e.g. 323.35 -> 323.35000000000002).
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "43.737.43"
$ws.Range("E2").Value = "  +1.38%  "
$ws.Range("D3").Value = "2.247.90"
$ws.Range("E3").Value = "  +0.82%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "323.35"
$ws.Range("E5").Value = "  +2.33%  "
$ws.Range("D6").Value = "101.96"
$ws.Range("E6").Value = "  +0.79%  "
$ws.Range("D7").Value = "0.581"
$ws.Range("E7").Value = "  -1.03%  "
$ws.Range("D9").Value = "0.557"
$ws.Range("E9").Value = "  -0.64%  "
$ws.Range("D10").Value = "37.43"
$ws.Range("E10").Value = "  +1.39%  "
$ws.Range("D11").Value = "0.0831"
$ws.Range("E11").Value = "  +0.95%  "
$ws.Range("D12").Value = "7.73"
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("E13").Value = "  -1.93%  "
$ws.Range("D14").Value = "2.587.82"
$ws.Range("D15").Value = "0.861"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "14.26"
$ws.Range("E16").Value = "  -0.50%  "
$ws.Range("D17").Value = "2.249.21"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "43.655.74"
$ws.Range("E18").Value = "  +1.28%  "
$ws.Range("D19").Value = "13.80"
$ws.Range("E19").Value = "  -3.56%  "
$ws.Range("D20").Value = "0.0₃0988"
$ws.Range("E20").Value = "  +2.90%  "
$ws.Range("D21").Value = "6.62"
$ws.Range("E21").Value = "  +2.06%  "
$ws.Range("D22").Value = "65.34"
$ws.Range("E22").Value = "  -0.23%  "
$ws.Range("D23").Value = "3.17"
$ws.Range("E23").Value = "  -0.72%  "
$ws.Range("D24").Value = "236.84"
$ws.Range("E24").Value = "  -0.44%  "
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").Value = "  +1.75%  "
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "10.15"
$ws.Range("E27").Value = "  +1.66%  "
$ws.Range("E28").Value = "  -1.33%  "
$ws.Range("D29").Value = "37.13"
$ws.Range("E29").Value = "  +7.70%  "
$ws.Range("D30").Value = "6.30"
$ws.Range("E30").Value = "  -1.25%  "
$ws.Range("D31").Value = "160.24"
$ws.Range("E31").Value = "  +3.88%  "
$ws.Range("D32").Value = "20.23"
$ws.Range("E32").Value = "  -1.23%  "
$ws.Range("D33").Value = "0.0856"
$ws.Range("E33").Value = "  -1.61%  "
$ws.Range("E34").Value = "  -2.63%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("E36").Value = "  +9.35%  "
$ws.Range("D37").Value = "1.93"
$ws.Range("E37").Value = "  +0.08%  "
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("E39").Value = "  +2.60%  "
$ws.Range("E40").Value = "  -3.46%  "
$ws.Range("D41").Value = "15.90"
$ws.Range("E41").Value = "  +20.55%  "
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("D44").Value = "1.811.08"
$ws.Range("E44").Value = "  +0.56%  "
$ws.Range("E45").Value = "  -2.19%  "
$ws.Range("D46").Value = "82.50"
$ws.Range("E46").Value = "  -6.04%  "
$ws.Range("D47").Value = "1.72"
$ws.Range("E47").Value = "  +6.81%  "
$ws.Range("D48").Value = "5.22"
$ws.Range("E48").Value = "  -2.04%  "
$ws.Range("D49").Value = "74.61"
$ws.Range("E49").Value = "  -3.03%  "
$ws.Range("D50").Value = "58.96"
$ws.Range("E50").Value = "  -0.11%  "
$ws.Range("D51").Value = "103.66"
$ws.Range("E51").Value = "  +0.35%  "
